# Fill in the March 2024 statistics, which were previously blank.
# (The "Yearly totals" sheet sums January..December via formulas, so it
# recalculates automatically once March's figures are entered here.)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("March")

$ws.Range("B2").Value = 1358
$ws.Range("C2").Value = 1200
$ws.Range("D2").Value = 158
$ws.Range("E2").Value = "We borrowerd more than we lent"
$ws.Range("G2").Value = "1.13 : 1"

$ws.Range("B3").Value = 619
$ws.Range("C3").Value = 487
$ws.Range("D3").Value = 132
$ws.Range("E3").Value = "We borrowerd more than we lent"
$ws.Range("G3").Value = "1.27 : 1"

$ws.Range("B4").Value = 1290
$ws.Range("C4").Value = 1267
$ws.Range("D4").Value = 23
$ws.Range("E4").Value = "We borrowerd more than we lent"
$ws.Range("G4").Value = "1.02 : 1"

$ws.Range("B5").Value = 26
$ws.Range("C5").Value = 116
$ws.Range("D5").Value = -90
$ws.Range("F5").Value = "We lent more than we borrowed"
$ws.Range("G5").Value = "0.22 : 1"

$ws.Range("B6").Value = 1206
$ws.Range("C6").Value = 1584
$ws.Range("D6").Value = -378
$ws.Range("F6").Value = "We lent more than we borrowed"
$ws.Range("G6").Value = "0.76 : 1"

$ws.Range("B7").Value = 291
$ws.Range("C7").Value = 160
$ws.Range("D7").Value = 131
$ws.Range("E7").Value = "We borrowerd more than we lent"
$ws.Range("G7").Value = "1.82 : 1"

$ws.Range("B8").Value = 109
$ws.Range("C8").Value = 194
$ws.Range("D8").Value = -85
$ws.Range("F8").Value = "We lent more than we borrowed"
$ws.Range("G8").Value = "0.56 : 1"

$ws.Range("B9").Value = 38
$ws.Range("C9").Value = 76
$ws.Range("D9").Value = -38
$ws.Range("F9").Value = "We lent more than we borrowed"
$ws.Range("G9").Value = "0.50 : 1"

$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 68
$ws.Range("D10").Value = -68
$ws.Range("F10").Value = "We lent more than we borrowed"
$ws.Range("G10").Value = "0.00 : 1"

$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0

$ws.Range("B12").Value = 21
$ws.Range("C12").Value = 26
$ws.Range("D12").Value = -5
$ws.Range("F12").Value = "We lent more than we borrowed"
$ws.Range("G12").Value = "0.81 : 1"

$ws.Range("B13").Value = 166
$ws.Range("C13").Value = 96
$ws.Range("D13").Value = 70
$ws.Range("E13").Value = "We borrowerd more than we lent"
$ws.Range("G13").Value = "1.73 : 1"

$ws.Range("B14").Value = 104
$ws.Range("C14").Value = 270
$ws.Range("D14").Value = -166
$ws.Range("F14").Value = "We lent more than we borrowed"
$ws.Range("G14").Value = "0.39 : 1"

$ws.Range("B15").Value = 44
$ws.Range("C15").Value = 133
$ws.Range("D15").Value = -89
$ws.Range("F15").Value = "We lent more than we borrowed"
$ws.Range("G15").Value = "0.33 : 1"

$ws.Range("B16").Value = 70
$ws.Range("C16").Value = 137
$ws.Range("D16").Value = -67
$ws.Range("F16").Value = "We lent more than we borrowed"
$ws.Range("G16").Value = "0.51 : 1"

$ws.Range("B17").Value = 610
$ws.Range("C17").Value = 386
$ws.Range("D17").Value = 224
$ws.Range("E17").Value = "We borrowerd more than we lent"
$ws.Range("G17").Value = "1.58 : 1"

$ws.Range("B18").Value = 64
$ws.Range("C18").Value = 82
$ws.Range("D18").Value = -18
$ws.Range("F18").Value = "We lent more than we borrowed"
$ws.Range("G18").Value = "0.78 : 1"

$ws.Range("B19").Value = 555
$ws.Range("C19").Value = 543
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = "We borrowerd more than we lent"
$ws.Range("G19").Value = "1.02 : 1"

$ws.Range("B20").Value = 2
$ws.Range("C20").Value = 94
$ws.Range("D20").Value = -92
$ws.Range("F20").Value = "We lent more than we borrowed"
$ws.Range("G20").Value = "0.02 : 1"

$ws.Range("B21").Value = 464
$ws.Range("C21").Value = 338
$ws.Range("D21").Value = 126
$ws.Range("E21").Value = "We borrowerd more than we lent"
$ws.Range("G21").Value = "1.37 : 1"

$ws.Range("B22").Value = 27
$ws.Range("C22").Value = 126
$ws.Range("D22").Value = -99
$ws.Range("F22").Value = "We lent more than we borrowed"
$ws.Range("G22").Value = "0.21 : 1"

$ws.Range("B23").Value = 636
$ws.Range("C23").Value = 384
$ws.Range("D23").Value = 252
$ws.Range("E23").Value = "We borrowerd more than we lent"
$ws.Range("G23").Value = "1.66 : 1"

$ws.Range("B24").Value = 2112
$ws.Range("C24").Value = 1209
$ws.Range("D24").Value = 903
$ws.Range("E24").Value = "We borrowerd more than we lent"
$ws.Range("G24").Value = "1.75 : 1"

$ws.Range("B25").Value = 168
$ws.Range("C25").Value = 324
$ws.Range("D25").Value = -156
$ws.Range("F25").Value = "We lent more than we borrowed"
$ws.Range("G25").Value = "0.52 : 1"

$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0

$ws.Range("B27").Value = 264
$ws.Range("C27").Value = 171
$ws.Range("D27").Value = 93
$ws.Range("E27").Value = "We borrowerd more than we lent"
$ws.Range("G27").Value = "1.54 : 1"

$ws.Range("B28").Value = 73
$ws.Range("C28").Value = 81
$ws.Range("D28").Value = -8
$ws.Range("F28").Value = "We lent more than we borrowed"
$ws.Range("G28").Value = "0.90 : 1"

$ws.Range("B29").Value = 662
$ws.Range("C29").Value = 418
$ws.Range("D29").Value = 244
$ws.Range("E29").Value = "We borrowerd more than we lent"
$ws.Range("G29").Value = "1.58 : 1"

$ws.Range("B30").Value = 47
$ws.Range("C30").Value = 22
$ws.Range("D30").Value = 25
$ws.Range("E30").Value = "We borrowerd more than we lent"
$ws.Range("G30").Value = "2.14 : 1"

$ws.Range("B31").Value = 52
$ws.Range("C31").Value = 341
$ws.Range("D31").Value = -289
$ws.Range("F31").Value = "We lent more than we borrowed"
$ws.Range("G31").Value = "0.15 : 1"

$ws.Range("B32").Value = 412
$ws.Range("C32").Value = 584
$ws.Range("D32").Value = -172
$ws.Range("F32").Value = "We lent more than we borrowed"
$ws.Range("G32").Value = "0.71 : 1"

$ws.Range("B33").Value = 299
$ws.Range("C33").Value = 632
$ws.Range("D33").Value = -333
$ws.Range("F33").Value = "We lent more than we borrowed"
$ws.Range("G33").Value = "0.47 : 1"

$ws.Range("B34").Value = 182
$ws.Range("C34").Value = 93
$ws.Range("D34").Value = 89
$ws.Range("E34").Value = "We borrowerd more than we lent"
$ws.Range("G34").Value = "1.96 : 1"

$ws.Range("B35").Value = 865
$ws.Range("C35").Value = 1210
$ws.Range("D35").Value = -345
$ws.Range("F35").Value = "We lent more than we borrowed"
$ws.Range("G35").Value = "0.71 : 1"

$ws.Range("B36").Value = 220
$ws.Range("C36").Value = 444
$ws.Range("D36").Value = -224
$ws.Range("F36").Value = "We lent more than we borrowed"
$ws.Range("G36").Value = "0.50 : 1"

$ws.Range("B37").Value = 487
$ws.Range("C37").Value = 322
$ws.Range("D37").Value = 165
$ws.Range("E37").Value = "We borrowerd more than we lent"
$ws.Range("G37").Value = "1.51 : 1"

$ws.Range("B38").Value = 27
$ws.Range("C38").Value = 158
$ws.Range("D38").Value = -131
$ws.Range("F38").Value = "We lent more than we borrowed"
$ws.Range("G38").Value = "0.17 : 1"

$ws.Range("B39").Value = 23
$ws.Range("C39").Value = 91
$ws.Range("D39").Value = -68
$ws.Range("F39").Value = "We lent more than we borrowed"
$ws.Range("G39").Value = "0.25 : 1"

$ws.Range("B40").Value = 62
$ws.Range("C40").Value = 103
$ws.Range("D40").Value = -41
$ws.Range("F40").Value = "We lent more than we borrowed"
$ws.Range("G40").Value = "0.60 : 1"

$ws.Range("B41").Value = 3
$ws.Range("C41").Value = 37
$ws.Range("D41").Value = -34
$ws.Range("F41").Value = "We lent more than we borrowed"
$ws.Range("G41").Value = "0.08 : 1"

$ws.Range("B42").Value = 9
$ws.Range("C42").Value = 22
$ws.Range("D42").Value = -13
$ws.Range("F42").Value = "We lent more than we borrowed"
$ws.Range("G42").Value = "0.41 : 1"

$ws.Range("B43").Value = 0
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 0

$ws.Range("B44").Value = 109
$ws.Range("C44").Value = 74
$ws.Range("D44").Value = 35
$ws.Range("E44").Value = "We borrowerd more than we lent"
$ws.Range("G44").Value = "1.47 : 1"

$ws.Range("B45").Value = 85
$ws.Range("C45").Value = 207
$ws.Range("D45").Value = -122
$ws.Range("F45").Value = "We lent more than we borrowed"
$ws.Range("G45").Value = "0.41 : 1"

$ws.Range("B46").Value = 444
$ws.Range("C46").Value = 676
$ws.Range("D46").Value = -232
$ws.Range("F46").Value = "We lent more than we borrowed"
$ws.Range("G46").Value = "0.66 : 1"

$ws.Range("B47").Value = 1079
$ws.Range("C47").Value = 569
$ws.Range("D47").Value = 510
$ws.Range("E47").Value = "We borrowerd more than we lent"
$ws.Range("G47").Value = "1.90 : 1"

$ws.Range("B48").Value = 299
$ws.Range("C48").Value = 580
$ws.Range("D48").Value = -281
$ws.Range("F48").Value = "We lent more than we borrowed"
$ws.Range("G48").Value = "0.52 : 1"

$ws.Range("B49").Value = 585
$ws.Range("C49").Value = 232
$ws.Range("D49").Value = 353
$ws.Range("E49").Value = "We borrowerd more than we lent"
$ws.Range("G49").Value = "2.52 : 1"

$ws.Range("B50").Value = 924
$ws.Range("C50").Value = 552
$ws.Range("D50").Value = 372
$ws.Range("E50").Value = "We borrowerd more than we lent"
$ws.Range("G50").Value = "1.67 : 1"

$ws.Range("B51").Value = 271
$ws.Range("C51").Value = 186
$ws.Range("D51").Value = 85
$ws.Range("E51").Value = "We borrowerd more than we lent"
$ws.Range("G51").Value = "1.46 : 1"

$ws.Range("B52").Value = 378
$ws.Range("C52").Value = 477
$ws.Range("D52").Value = -99
$ws.Range("F52").Value = "We lent more than we borrowed"
$ws.Range("G52").Value = "0.79 : 1"

$ws.Range("B53").Value = 127
$ws.Range("C53").Value = 229
$ws.Range("D53").Value = -102
$ws.Range("F53").Value = "We lent more than we borrowed"
$ws.Range("G53").Value = "0.55 : 1"

$ws.Range("B54").Value = 14
$ws.Range("C54").Value = 201
$ws.Range("D54").Value = -187
$ws.Range("F54").Value = "We lent more than we borrowed"
$ws.Range("G54").Value = "0.07 : 1"

$ws.Range("B55").Value = 268
$ws.Range("C55").Value = 238
$ws.Range("D55").Value = 30
$ws.Range("E55").Value = "We borrowerd more than we lent"
$ws.Range("G55").Value = "1.13 : 1"
